# basic_evaluation.xlsx - "using named_range in a cell without any operation works correctly"
#
# 1. A17 moves from 10 to 30 (the rest of row 17 are formulas that key off of
#    it/Liste3, so their cached results ripple automatically on recalculation).
# 2. New column L holds a plain `=Liste2` formula (a named range used bare, with
#    no surrounding operation) in L1:L3, resolving via implicit intersection to
#    the matching row of B1:B3 (10, 20, 30).
# 3. The sheet selection is left on the new L1:L3 block, active cell L1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17").Value = 30

$ws.Range("L1").Formula = "=Liste2"
$ws.Range("L2").Formula = "=Liste2"
$ws.Range("L3").Formula = "=Liste2"

$ws.Range("L1:L3").Select() | Out-Null
